# Updated CHE_grids model - 2025-08-09 17:35
# Reassign the "grid_cell" (AG column) values on the "solar" sheet's
# distr_elc_won-CHE_xxxx table (rows 4-26) to their new grid-cell mapping.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("solar")

$ws.Range("AG4").Value2  = "CHE_11"
$ws.Range("AG5").Value2  = "CHE_15"
$ws.Range("AG6").Value2  = "CHE_25"
$ws.Range("AG7").Value2  = "CHE_24"
$ws.Range("AG8").Value2  = "CHE_5"
$ws.Range("AG9").Value2  = "CHE_8"
$ws.Range("AG10").Value2 = "CHE_14"
$ws.Range("AG11").Value2 = "CHE_18"
$ws.Range("AG12").Value2 = "CHE_7"
$ws.Range("AG13").Value2 = "CHE_13"
$ws.Range("AG14").Value2 = "CHE_21"
$ws.Range("AG15").Value2 = "CHE_9"
$ws.Range("AG16").Value2 = "CHE_4"
$ws.Range("AG17").Value2 = "CHE_20"
$ws.Range("AG18").Value2 = "CHE_1"
$ws.Range("AG19").Value2 = "CHE_6"
$ws.Range("AG20").Value2 = "CHE_0"
$ws.Range("AG21").Value2 = "CHE_3"
$ws.Range("AG22").Value2 = "CHE_10"
$ws.Range("AG23").Value2 = "CHE_22"
$ws.Range("AG24").Value2 = "CHE_12"
$ws.Range("AG25").Value2 = "CHE_17"
$ws.Range("AG26").Value2 = "CHE_19"

$wb.Save()
